$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 71; this shifts the former rows 71-81 down to 72-82,
# carrying their formatting (e.g. the date style on column D) along with them.
$ws.Rows.Item(71).Insert()

# Populate the newly inserted row 71 with the new weekly record.
$ws.Cells.Item(71, 1).Value = 10
$ws.Cells.Item(71, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(71, 3).Value = "La Araucanía"
$ws.Cells.Item(71, 4).Value = 44491
$ws.Cells.Item(71, 5).Value = 9
$ws.Cells.Item(71, 6).Value = "Fruta"
$ws.Cells.Item(71, 7).Value = 100107
$ws.Cells.Item(71, 8).Value = "Otros"
$ws.Cells.Item(71, 9).Value = 100107002
$ws.Cells.Item(71, 10).Value = "Chirimoya"
$ws.Cells.Item(71, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(71, 12).Value = "Primera"
$ws.Cells.Item(71, 13).Value = 80
$ws.Cells.Item(71, 14).Value = 3000
$ws.Cells.Item(71, 15).Value = 3000
$ws.Cells.Item(71, 16).Value = 3000
$ws.Cells.Item(71, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(71, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(71, 19).Value = 3000
$ws.Cells.Item(71, 20).Value = 1
